$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.075.43"
$ws.Range("E2").Value = "  +0.54%  "
$ws.Range("D3").Value = "2.298.66"
$ws.Range("E3").Value = "  +0.32%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.00"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.50"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.44%  "
$ws.Range("E7").Value = "  +3.00%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.514"
$ws.Range("D9").ClearFormats()
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.17"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.52%  "
$ws.Range("E11").Value = "  +0.31%  "
$ws.Range("E12").Value = "  +0.58%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "17.73"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -2.90%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.87"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.70%  "
$ws.Range("D15").Value = "2.656.00"
$ws.Range("E15").Value = "  +0.23%  "
$ws.Range("D16").Value = "2.300.60"
$ws.Range("E16").Value = "  +0.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.787"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.40%  "
$ws.Range("D18").Value = "42.944.06"
$ws.Range("E18").Value = "  +0.44%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.04"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +4.97%  "
$ws.Range("D20").Value = "0.0₃0909"
$ws.Range("E20").Value = "  +0.92%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.12"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.17"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.77%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.58"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.86%  "
$ws.Range("E24").Value = "  -1.83%  "
$ws.Range("E25").Value = "  -0.36%  "
$ws.Range("E26").Value = "  -0.64%  "
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.98"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.26%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.04"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -13.14%  "
$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.15"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.29%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "162.92"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.54%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "32.97"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -4.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.07%  "
$ws.Range("E34").Value = "  +2.44%  "
$ws.Range("E35").Value = "  +2.85%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.72"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.60%  "
$ws.Range("E37").Value = "  +0.34%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0693"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.05%  "
$ws.Range("E39").Value = "  +1.09%  "
$ws.Range("E40").Value = "  -0.36%  "
$ws.Range("E41").Value = "  +1.36%  "
$ws.Range("E42").Value = "  -2.21%  "
$ws.Range("D43").Value = "2.009.26"
$ws.Range("E43").Value = "  +2.04%  "
$ws.Range("E44").Value = "  -1.53%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.19"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -5.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.23"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.86%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "17.50"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.02%  "
$ws.Range("E48").Value = "  -1.62%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "54.27"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.67%  "
$ws.Range("D50").Value = "2.529.16"
$ws.Range("E50").Value = "  +0.39%  "
$ws.Range("E51").Value = "  -0.44%  "
